$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The test data got re-sorted; swap the data that used to live in row 7
# into row 8, and vice versa (cylinder/engine/N/mean/sd columns E:K).
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 83.33333333333333
$ws.Range("I7").Value = 18.50225211517056
$ws.Range("J7").Value = 2.886666666666667
$ws.Range("K7").Value = 0.4911551010967242

$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 110
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2.7475
$ws.Range("K8").Value = 0.1803122292025695

# Row 9 keeps its own F:K stats, but its cylinder count (E9) now duplicates
# row 8's (both are "6"), so the two cylinder cells get merged into one,
# vertically centred-at-top, and E9 itself goes blank.
$ws.Range("E8:E9").Merge()

$ws.Range("L9").Copy()
$ws.Range("E9").PasteSpecial(-4122)

$ws.Range("E8").VerticalAlignment = -4160
